$wb = $excel.ActiveWorkbook

# Get the "місто Київ" sheet (last sheet, index 24) and rename it to "м. Київ"
$ws = $wb.Worksheets.Item("місто Київ")
$ws.Name = "м. Київ"

# Update the "Область" (R) column values on rows 2-9 from "місто Київ" to "м. Київ"
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 18)
    if ($cell.Value2 -eq "місто Київ") {
        $cell.Value = "м. Київ"
    }
}
